$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1: proper column headers (previously duplicated row-2 data) ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# header cells share the bold/centered/bordered style used by B1:G1
$ws.Range("H1:N1").Font.Bold = $true
$ws.Range("H1:N1").HorizontalAlignment = -4108
$ws.Range("H1:N1").VerticalAlignment = -4160
$ws.Range("H1:N1").Borders.LineStyle = 1

# --- Row 2: extend the car record with the remaining metadata columns ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# avoid Excel's automatic "looks like a date" literal -> serial-number coercion
# for the ISO-formatted date column, so it stays a plain text value
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2011-11-17"
$ws.Range("J2").Style = $ws.Range("G2").Style

$ws.Range("K2").Value = "蘇震清"
$ws.Range("L2").Value = 1718
$ws.Range("M2").Value = "tmp98701"
$ws.Range("N2").Value = 39
